# Update crypto price/volume figures per the Feb 16 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="51.881.73"; E="  +0.21%  "},
    @{Row=3; D="2.780.06"; E="  -1.86%  "},
    @{Row=4; E="  -0.05%  "},
    @{Row=5; D="357.43"; E="  +1.38%  "},
    @{Row=6; D="109.33"; E="  -3.57%  "},
    @{Row=7; E="  +0.91%  "},
    @{Row=8; D="0.999"; E="  -0.01%  "},
    @{Row=9; E="  -0.69%  "},
    @{Row=10; D="40.03"; E="  -3.69%  "},
    @{Row=11; E="  +0.07%  "},
    @{Row=12; E="  +0.75%  "},
    @{Row=13; D="19.45"; E="  -2.37%  "},
    @{Row=14; E="  -1.50%  "},
    @{Row=15; D="3.213.66"; E="  -2.04%  "},
    @{Row=16; D="2.794.15"; E="  -1.07%  "},
    @{Row=17; E="  +4.36%  "},
    @{Row=18; D="51.786.44"; E="  +0.18%  "},
    @{Row=19; E="  +0.73%  "},
    @{Row=20; E="  -0.23%  "},
    @{Row=21; D="13.03"; E="  -3.18%  "},
    @{Row=22; D="0.0₃0978"; E="  -1.46%  "},
    @{Row=23; D="273.68"; E="  +1.12%  "},
    @{Row=24; D="70.05"; E="  +0.52%  "},
    @{Row=25; E="  -1.12%  "},
    @{Row=26; D="26.65"; E="  -0.13%  "},
    @{Row=27; E="  -0.03%  "},
    @{Row=28; D="10.16"; E="  -1.08%  "},
    @{Row=29; E="  +3.91%  "},
    @{Row=30; E="  -1.40%  "},
    @{Row=31; D="0.0465"; E="  +4.15%  "},
    @{Row=32; D="51.56"; E="  +1.77%  "},
    @{Row=33; D="33.94"; E="  +0.09%  "},
    @{Row=34; E="  -1.75%  "},
    @{Row=35; D="0.0846"; E="  +2.45%  "},
    @{Row=36; E="  +7.84%  "},
    @{Row=37; E="  -0.02%  "},
    @{Row=38; D="3.25"; E="  +1.20%  "},
    @{Row=39; D="18.10"; E="  +0.58%  "},
    @{Row=40; E="  -3.67%  "},
    @{Row=41; E="  -0.35%  "},
    @{Row=42; D="2.53"; E="  -1.02%  "},
    @{Row=43; D="121.93"; E="  -2.91%  "},
    @{Row=44; E="  -2.60%  "},
    @{Row=45; D="22.13"; E="  -6.57%  "},
    @{Row=46; D="2.067.28"; E="  -0.59%  "},
    @{Row=47; E="  -2.31%  "},
    @{Row=48; E="  -5.86%  "},
    @{Row=49; E="  -0.13%  "},
    @{Row=50; D="0.935"; E="  +0.03%  "},
    @{Row=51; E="  +0.29%  "}
)

foreach ($item in $updates) {
    if ($item.ContainsKey("D")) {
        $cell = $ws.Cells.Item($item.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
        $cell.Style = "Normal"
    }
    if ($item.ContainsKey("E")) {
        $ws.Cells.Item($item.Row, 5).Value = $item.E
    }
}
